$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace "Davanagere" with "Davangere" in column G, rows 3..63,
# except where the cell doesn't literally equal "Davanagere"
# (row 18 and 59 have different school-name text, row 49 already
# reads "Davanager" and must stay untouched).
for ($r = 3; $r -le 63; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # column G
    if ($cell.Value2 -eq "Davanagere") {
        $cell.Value = "Davangere"
    }
}

# Clear the stray empty inline-string cells in column F for rows 18 and 59
# so they no longer exist (used range shrinks from H to G).
$ws.Cells.Item(18, 6).ClearContents()
$ws.Cells.Item(59, 6).ClearContents()
